# Apply "permutate DSM matrix" edit: reorder the 23 function rows/columns of the
# DSM matrix according to the new (permutated) order, carrying the 0/1 relation
# data and the diagonal "F#" function-id labels along by name. The
# LogicalComponentName row (row 25) is regenerated independently of the
# permutation (it is produced by the initial/default export path), so it is
# written from its own fixed sequence rather than being re-derived from the
# reordered columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of the 23 functions (used for both the column headers in row 1
# and the row headers in column A - the DSM is symmetric).
$names = @(
    "Send aircraft view",
    "Manage Mission Modes",
    "Sense and Avoid Obstacles",
    "Build FlightPlan Relative to Aircraft Type",
    "CheckWinfForce",
    "Retrieve POI",
    "Identify Absolute Aircraft Coordinates",
    "Identify Defects",
    "Send Pictures to DB",
    "Generate Thrust",
    "Monitor UAV Control",
    "Manage Photos Recording",
    "Record photos and videos",
    "Control Camera Orientation",
    "Configurate Flight Plan",
    "Send moving obstacle position",
    "Send stationary obstacle position",
    "Send/Receive data",
    "Emergency Landing",
    "Control UAV Position",
    "Control UAV attitude",
    "Send perceived position, attitude ",
    "Send command and position setting"
)

# Permutated relation matrix (rows/cols 2..24 => function positions 1..23).
# Diagonal cells keep the worksheet's positional "F<n>" function-id label.
$grid = @(
    @("F1", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, "F2", 0, 0, 0, 1, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0),
    @(0, 0, "F3", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 1, 0, 0),
    @(0, 0, 0, "F4", 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 1, 0, 0, "F5", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, "F6", 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 1, 0, 0),
    @(0, 0, 0, 0, 0, 1, "F7", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, "F8", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 1, "F9", 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, "F10", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F11", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1),
    @(0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, "F12", 1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F13", 1, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F14", 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F15", 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F16", 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F17", 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F18", 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F19", 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F20", 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F21", 0, 0),
    @(0, 0, 0, 0, 1, 0, 1, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F22", 0),
    @(0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "F23")
)

# Row 25 (LogicalComponentName) - independent of the permutation above.
$row25 = @(
    "Aircraft",
    "Mission Mgt Subsystem",
    "Mission Mgt Subsystem",
    "Mission Mgt Subsystem",
    "Mission Mgt Subsystem",
    "Mission Mgt Subsystem",
    "Mission Mgt Subsystem",
    "Mission Mgt Subsystem",
    "Mission Mgt Subsystem",
    "Propulsion Subsystem",
    "UAV Control Station Subsystem",
    "Vision Subsystem",
    "Vision Subsystem",
    "Vision Subsystem",
    "Airline Human Operator",
    "Moving Obstacles",
    "Stationary Obstacle",
    "Aircraft Company Database",
    "Flight Mgt Subsystem",
    "Flight Mgt Subsystem",
    "Flight Mgt Subsystem",
    "Air/ Terrestrian Gravity",
    "UAV Pilot"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(1, $col).Value = $names[$i]
    $ws.Cells.Item($col, 1).Value = $names[$i]
    $ws.Cells.Item(25, $col).Value = $row25[$i]
}

for ($r = 0; $r -lt $grid.Length; $r++) {
    $row = $r + 2
    $rowData = $grid[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 2
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}

